# Scheduled runner refresh: update market-price / profit figures (columns
# H..N) on a handful of rows across the ALC / ARM / BSM / CRP / CUL / GSM /
# LTW / WVR leve-profit sheets, per the latest FFXIV market data pull.
# Cells whose computed profit is no longer applicable are cleared; cells
# that newly have an applicable profit are populated.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (G id 4564)
$ws.Range("H6").Value = 233.33333
$ws.Range("I6").Value = 233.33333
$ws.Range("K6").Value = 699.99999
$ws.Range("M6").Value = -587.99999

# Row 26 (G id 1963)
$ws.Range("H26").Value = 40000
$ws.Range("J26").Value = 40000
$ws.Range("L26").Value = 40000
$ws.Range("N26").Value = -40688

# Row 45 (G id 4585)
$ws.Range("H45").Value = 117
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Row 74 (G id 5507)
$ws.Range("H74").Value = 3566.6667
$ws.Range("J74").Value = 3566.6667
$ws.Range("L74").Value = 3566.6667
$ws.Range("N74").Value = -5438.6667

# Row 77 (G id 5507)
$ws.Range("H77").Value = 3566.6667
$ws.Range("J77").Value = 3566.6667
$ws.Range("L77").Value = 17833.3335
$ws.Range("N77").Value = -27193.3335

# Row 106 (G id 19903)
$ws.Range("H106").Value = 3710.7144
$ws.Range("I106").Value = 2661.5
$ws.Range("K106").Value = 2661.5
$ws.Range("M106").Value = -2030.5

# Row 111 (G id 27768)
$ws.Range("H111").Value = 2950.1052
$ws.Range("I111").Value = 3102.9
$ws.Range("K111").Value = 9308.700000000001
$ws.Range("M111").Value = -6241.700000000001

# Row 116 (G id 27778)
$ws.Range("H116").Value = 5778.857
$ws.Range("I116").Value = 2312.875
$ws.Range("J116").Value = 10400.167
$ws.Range("K116").Value = 2312.875
$ws.Range("L116").Value = 10400.167
$ws.Range("M116").Value = 1129.125
$ws.Range("N116").Value = -17284.167

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (G id 27713)
$ws.Range("H2").Value = 2301.25
$ws.Range("I2").Value = 2083
$ws.Range("J2").Value = 2665
$ws.Range("K2").Value = 2083
$ws.Range("L2").Value = 2665
$ws.Range("M2").Value = -1970
$ws.Range("N2").Value = -2891

# Row 45 (G id 27714)
$ws.Range("H45").Value = 1559.2
$ws.Range("I45").Value = 1226.6
$ws.Range("J45").Value = 1891.8
$ws.Range("K45").Value = 1226.6
$ws.Range("L45").Value = 1891.8
$ws.Range("M45").Value = -849.5999999999999
$ws.Range("N45").Value = -2645.8

# Row 74 (G id 44000)
$ws.Range("H74").Value = 56155.7
$ws.Range("I74").Value = 67826.734
$ws.Range("J74").Value = 21142.6
$ws.Range("K74").Value = 67826.734
$ws.Range("L74").Value = 21142.6
$ws.Range("M74").Value = -66952.734
$ws.Range("N74").Value = -22890.6

# Row 77 (G id 44000)
$ws.Range("H77").Value = 56155.7
$ws.Range("I77").Value = 67826.734
$ws.Range("J77").Value = 21142.6
$ws.Range("K77").Value = 339133.67
$ws.Range("L77").Value = 105713
$ws.Range("M77").Value = -334765.67
$ws.Range("N77").Value = -114449

# Row 116 (G id 27713)
$ws.Range("H116").Value = 2301.25
$ws.Range("I116").Value = 2083
$ws.Range("J116").Value = 2665
$ws.Range("K116").Value = 2083
$ws.Range("L116").Value = 2665
$ws.Range("M116").Value = 211
$ws.Range("N116").Value = -7253

# Row 125 (G id 34251)
$ws.Range("H125").Value = 19888.889
$ws.Range("J125").Value = 19888.889
$ws.Range("L125").Value = 19888.889
$ws.Range("N125").Value = -29728.889

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (G id 27713)
$ws.Range("H3").Value = 2301.25
$ws.Range("I3").Value = 2083
$ws.Range("J3").Value = 2665
$ws.Range("K3").Value = 2083
$ws.Range("L3").Value = 2665
$ws.Range("M3").Value = -1969
$ws.Range("N3").Value = -2893

# Row 32 (G id 2397)
$ws.Range("H32").Value = 27500
$ws.Range("J32").Value = 27500
$ws.Range("L32").Value = 27500
$ws.Range("N32").Value = -28268

# Row 33 (G id 1625)
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 21 (G id 2000)
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 31 (G id 44023)
$ws.Range("H31").Value = 30304726
$ws.Range("J31").Value = 3225.5
$ws.Range("L31").Value = 3225.5
$ws.Range("N31").Value = -3815.5

# Row 34 (G id 44023)
$ws.Range("H34").Value = 30304726
$ws.Range("J34").Value = 3225.5
$ws.Range("L34").Value = 3225.5
$ws.Range("N34").Value = -3629.5

# Row 62 (G id 12580)
$ws.Range("H62").Value = 3136.3635
$ws.Range("I62").Value = 2350
$ws.Range("J62").Value = 4080
$ws.Range("K62").Value = 2350
$ws.Range("L62").Value = 4080
$ws.Range("M62").Value = -1726
$ws.Range("N62").Value = -5328

# Row 65 (G id 12580)
$ws.Range("H65").Value = 3136.3635
$ws.Range("I65").Value = 2350
$ws.Range("J65").Value = 4080
$ws.Range("K65").Value = 11750
$ws.Range("L65").Value = 20400
$ws.Range("M65").Value = -8630
$ws.Range("N65").Value = -26640

# Row 75 (G id 11936)
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31996

# Row 78 (G id 11936)
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99984

$ws = $wb.Worksheets.Item("CUL")
# Row 70 (G id 12867)
$ws.Range("H70").Value = 5503.1816
$ws.Range("I70").Value = 3259
$ws.Range("J70").Value = 6785.5713
$ws.Range("K70").Value = 9777
$ws.Range("L70").Value = 20356.7139
$ws.Range("M70").Value = -9462
$ws.Range("N70").Value = -20986.7139

# Row 73 (G id 12867)
$ws.Range("H73").Value = 5503.1816
$ws.Range("I73").Value = 3259
$ws.Range("J73").Value = 6785.5713
$ws.Range("K73").Value = 9777
$ws.Range("L73").Value = 20356.7139
$ws.Range("M73").Value = -8685
$ws.Range("N73").Value = -22540.7139

# Row 129 (G id 36054)
$ws.Range("H129").Value = 3366.6667
$ws.Range("I129").Value = 3459.4
$ws.Range("J129").Value = 3312.1177
$ws.Range("K129").Value = 10378.2
$ws.Range("L129").Value = 9936.3531
$ws.Range("M129").Value = -5378.200000000001
$ws.Range("N129").Value = -19936.3531

# Row 133 (G id 44073)
$ws.Range("H133").Value = 7695.909
$ws.Range("I133").Value = 5327.5
$ws.Range("J133").Value = 8222.223
$ws.Range("K133").Value = 15982.5
$ws.Range("L133").Value = 24666.669
$ws.Range("M133").Value = -10922.5
$ws.Range("N133").Value = -34786.669

$ws = $wb.Worksheets.Item("GSM")
# Row 133 (G id 41854)
$ws.Range("H133").Value = 34543
$ws.Range("J133").Value = 34543
$ws.Range("L133").Value = 34543
$ws.Range("N133").Value = -44663

$ws = $wb.Worksheets.Item("LTW")
# Row 48 (G id 3625)
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

# Row 122 (G id 36247)
$ws.Range("H122").Value = 3826.2856
$ws.Range("I122").Value = 3642.5454
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 10927.6362
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -8477.636200000001
$ws.Range("N122").Value = -18400

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (G id 27746)
$ws.Range("H107").Value = 3902.1875
$ws.Range("I107").Value = 8056.7144
$ws.Range("J107").Value = 670.8889
$ws.Range("K107").Value = 24170.1432
$ws.Range("L107").Value = 2012.6667
$ws.Range("M107").Value = -22250.1432
$ws.Range("N107").Value = -5852.6667

# Row 122 (G id 36208)
$ws.Range("H122").Value = 108639.54
$ws.Range("I122").Value = 20730.4
$ws.Range("J122").Value = 401670
$ws.Range("K122").Value = 62191.2
$ws.Range("L122").Value = 1205010
$ws.Range("M122").Value = -59741.2
$ws.Range("N122").Value = -1209910
